# Updated Stargen for Elevated Plains
# Add a new "Higher Plains" section (rows 37-39) mirroring the existing
# "Plains" section (rows 33-35), but with new biome names in column B.
# Row 36, which used to hold the "Weblands / Flower Forest / hingeland but
# minecraft sussy" placeholder row, becomes the new section header
# ("    Higher Plains") and loses its B/C values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36 becomes the new "Higher Plains" section header; clear B36/C36.
$ws.Range("A36").Value = "    Higher Plains"
$ws.Range("B36").ClearContents()
$ws.Range("C36").ClearContents()

# Populate column B first, in the order that matches the shared-string
# insertion order used by the canonical workbook (Taiga Mountains, then
# Snowy Taiga, then Snowy Tundra), before filling in A/C so the shared
# string table indices line up.
$ws.Range("B39").Value = "Taiga Mountains"
$ws.Range("B38").Value = "Snowy Taiga"
$ws.Range("B37").Value = "Snowy Tundra"

# Row 37: mirrors row 33 (Flatlands / Regular Plains Biome) with a new B value.
$ws.Range("A37").Value = $ws.Range("A33").Value2
$ws.Range("C37").Value = $ws.Range("C33").Value2

# Row 38: mirrors row 34 (Forest / Better version of the Forest Biome).
$ws.Range("A38").Value = $ws.Range("A34").Value2
$ws.Range("C38").Value = $ws.Range("C34").Value2

# Row 39: mirrors row 35 (Hilled Plains / Bigger Plains).
$ws.Range("A39").Value = $ws.Range("A35").Value2
$ws.Range("C39").Value = $ws.Range("C35").Value2

# Update the view: scroll so row 15 is at the top, and select C37.
$excel.ActiveWindow.ScrollRow = 15
$ws.Range("C37").Select()
